# Update "Forecast Comparison" sheet with a new Week_Start_Date column,
# shorten the Week labels (W01 -> W1, ... W09 -> W9), and store
# is_holiday_week as a proper boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert a new column B ("Week_Start_Date"), pushing ASIN .. is_holiday_week
# --- one column to the right (B:I -> C:J). The sheet dimension grows to A1:J17.
$ws.Range("B1").EntireColumn.Insert()

# --- New header ---
$ws.Range("B1").Value = "Week_Start_Date"

# --- Week labels (column A) lose their leading zero for weeks 1-9 ---
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $weeks[$i]
}

# --- Week start dates (column B), rows 2-17, one per week. Stored as plain
# --- text (leading apostrophe keeps Excel from coercing them to date serials,
# --- matching the source inlineStr "YYYY-MM-DD" cells). ---
$startDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)
for ($i = 0; $i -lt $startDates.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = "'" + $startDates[$i]
}

# --- is_holiday_week (column J after the insert) becomes a boolean FALSE
# --- instead of the numeric 0 it inherited from the shift. ---
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}
